# Weekly CompStat refresh: roll the report forward one week and update the
# crime-complaint figures (new crime data collected).
#
# xlPasteFormats = -4122 (used to clone a cell's number-format / style onto
# another cell before changing its value, so cells that switch between a
# numeric display and the "placeholder" text display - shared strings
# "0" / "***.*" - end up with the right <c s="..."> style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Copy-Style($ws, $fromRef, $toRef) {
    $ws.Range($fromRef).Copy()
    $ws.Range($toRef).PasteSpecial(-4122)
}

# Set a plain numeric value, keeping the cell's current style untouched.
function Set-Num($ws, $cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

# Turn a cell into a numeric cell with a specific style (borrowed from
# $styleFromRef), e.g. converting a "0"/"***.*" placeholder back into a
# real number.
function Set-NumWithStyle($ws, $cellRef, $val, $styleFromRef) {
    Copy-Style $ws $styleFromRef $cellRef
    $ws.Range($cellRef).Value = $val
}

# Turn a cell into a text "placeholder" cell (shared string "0" or
# "***.*"), keeping the look of a numeric column by borrowing the style
# from $styleFromRef (typically a neighboring already-placeholder cell).
function Set-TextPlaceholder($ws, $cellRef, $text, $styleFromRef) {
    Copy-Style $ws $styleFromRef $cellRef
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    Copy-Style $ws $styleFromRef $cellRef
}

# ---------------------------------------------------------------------
# Header text: bump the issue number and roll the reporting week forward.
# ---------------------------------------------------------------------

# A8 = "Volume 32   Number  15" -> "...Number  16"
$cell = $ws.Range("A8")
$cell.Characters(21, 2).Text = "16"
$len = $cell.Characters().Text.Length
$cell.Characters(1, 20).Font.Size = 10
$cell.Characters(21, $len - 20).Font.Size = 10

# C9 = "Report Covering the Week  4/7/2025  Through  4/13/2025"
#   -> "Report Covering the Week  4/14/2025  Through  4/20/2025"
$cell = $ws.Range("C9")
$cell.Characters(27, 8).Text = "4/14/2025"
$cell.Characters(47, 9).Text = "4/20/2025"
$len = $cell.Characters().Text.Length
$cell.Characters(1, 26).Font.Size = 10
$cell.Characters(27, $len - 26).Font.Size = 10

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
Set-TextPlaceholder $ws "G14" "0" "F14"
Set-TextPlaceholder $ws "H14" "***.*" "F14"

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextPlaceholder $ws "C15" "0" "F14"
Set-TextPlaceholder $ws "D15" "0" "F14"
Set-TextPlaceholder $ws "E15" "***.*" "F14"
Set-Num $ws "L15" 150

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
Set-Num $ws "C16" 2
Set-NumWithStyle $ws "D16" 1 "C16"
Set-NumWithStyle $ws "E16" 100 "H16"
Set-Num $ws "F16" 8
Set-Num $ws "H16" 100
Set-Num $ws "I16" 18
Set-Num $ws "J16" 25
Set-Num $ws "K16" -28
Set-Num $ws "L16" -30.769230769230

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
Set-Num $ws "C17" 1
Set-Num $ws "D17" 9
Set-Num $ws "E17" -88.888888888888
Set-Num $ws "F17" 11
Set-Num $ws "G17" 18
Set-Num $ws "H17" -38.888888888888
Set-Num $ws "I17" 67
Set-Num $ws "J17" 70
Set-Num $ws "K17" -4.285714285714
Set-Num $ws "L17" 3.076923076923

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
Set-TextPlaceholder $ws "F18" "0" "C18"
Set-Num $ws "H18" -100
Set-Num $ws "J18" 18
Set-Num $ws "K18" -27.777777777777

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
Set-Num $ws "C19" 4
Set-Num $ws "D19" 5
Set-Num $ws "E19" -20
Set-Num $ws "F19" 20
Set-Num $ws "G19" 21
Set-Num $ws "H19" -4.761904761904
Set-Num $ws "I19" 70
Set-Num $ws "J19" 88
Set-Num $ws "K19" -20.454545454545
Set-Num $ws "L19" -15.662650602409

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-Num $ws "C20" 6
Set-Num $ws "E20" 200
Set-Num $ws "F20" 20
Set-Num $ws "G20" 12
Set-Num $ws "H20" 66.666666666666
Set-Num $ws "I20" 63
Set-Num $ws "J20" 47
Set-Num $ws "K20" 34.042553191489
Set-Num $ws "L20" 12.5

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
Set-Num $ws "C21" 13
Set-Num $ws "D21" 18
Set-Num $ws "E21" -27.777777777777
Set-Num $ws "F21" 61
Set-Num $ws "G21" 62
Set-Num $ws "H21" -1.612903225806
Set-Num $ws "I21" 241
Set-Num $ws "J21" 253
Set-Num $ws "K21" -4.743083003952
Set-Num $ws "L21" -4.743083003952

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
Set-Num $ws "C24" 5
Set-Num $ws "D24" 16
Set-Num $ws "E24" -68.75
Set-Num $ws "F24" 41
Set-Num $ws "G24" 46
Set-Num $ws "H24" -10.869565217391
Set-Num $ws "I24" 151
Set-Num $ws "J24" 187
Set-Num $ws "K24" -19.251336898395
Set-Num $ws "L24" -12.209302325581

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
Set-NumWithStyle $ws "D25" 4 "C25"
Set-NumWithStyle $ws "E25" -50 "H25"
Set-Num $ws "G25" 12
Set-Num $ws "H25" -33.333333333333
Set-Num $ws "I25" 35
Set-Num $ws "J25" 36
Set-Num $ws "K25" -2.777777777777
Set-Num $ws "L25" 16.666666666666

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
Set-Num $ws "C26" 3
Set-Num $ws "D26" 8
Set-Num $ws "E26" -62.5
Set-Num $ws "F26" 24
Set-Num $ws "G26" 19
Set-Num $ws "H26" 26.315789473684
Set-Num $ws "I26" 104
Set-Num $ws "J26" 121
Set-Num $ws "K26" -14.049586776859
Set-Num $ws "L26" 23.809523809523

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
Set-Num $ws "C27" 1
Set-TextPlaceholder $ws "D27" "0" "F14"
Set-TextPlaceholder $ws "E27" "***.*" "F14"
Set-Num $ws "F27" 3
Set-Num $ws "G27" 3
Set-Num $ws "H27" 0
Set-Num $ws "I27" 13
Set-Num $ws "K27" 116.666666666667
Set-Num $ws "L27" 30

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
Set-NumWithStyle $ws "D28" 2 "F28"
Set-NumWithStyle $ws "E28" -100 "K28"
Set-NumWithStyle $ws "G28" 2 "F28"
Set-NumWithStyle $ws "H28" -50 "K28"
Set-Num $ws "J28" 4
Set-Num $ws "K28" -25
